# Apply updated cryptos list values (price + 1h volume change)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to hold a literal text string (matches the
    # original inline-string cell type) instead of letting Excel
    # auto-coerce numeric-looking text ("235.70", "0.6023", ...)
    # into a real number. Revert the style afterwards so no
    # lingering NumberFormat is left on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.291.85"
$ws.Range("E2").Value = "  -0.12%  "
Set-TextValue $ws.Range("D3") "1.829.84"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.42%  "
Set-TextValue $ws.Range("D5") "235.70"
$ws.Range("E5").Value = "  -1.38%  "
Set-TextValue $ws.Range("D6") "0.6023"
$ws.Range("E6").Value = "  -3.18%  "
$ws.Range("E7").Value = "  +0.43%  "
Set-TextValue $ws.Range("D8") "0.06969"
$ws.Range("E8").Value = "  -4.90%  "
$ws.Range("E9").Value = "  -3.61%  "
Set-TextValue $ws.Range("D10") "23.57"
$ws.Range("E10").Value = "  -4.37%  "
Set-TextValue $ws.Range("D11") "0.07623"
$ws.Range("E11").Value = "  -1.36%  "
Set-TextValue $ws.Range("D12") "1.835.34"
$ws.Range("E12").Value = "  +0.68%  "
Set-TextValue $ws.Range("D13") "4.754"
$ws.Range("E13").Value = "  -3.75%  "
Set-TextValue $ws.Range("D14") "0.6322"
$ws.Range("E14").Value = "  -4.05%  "
Set-TextValue $ws.Range("D15") "0.000009852"
$ws.Range("E15").Value = "  -4.20%  "
Set-TextValue $ws.Range("D16") "77.80"
$ws.Range("E16").Value = "  -4.36%  "
Set-TextValue $ws.Range("D17") "28.986.98"
$ws.Range("E17").Value = "  -1.13%  "
Set-TextValue $ws.Range("D18") "5.589"
$ws.Range("E18").Value = "  -10.49%  "
Set-TextValue $ws.Range("D19") "217.64"
$ws.Range("E19").Value = "  -7.90%  "
$ws.Range("E20").Value = "  +0.35%  "
Set-TextValue $ws.Range("D21") "11.59"
$ws.Range("E21").Value = "  -4.84%  "
Set-TextValue $ws.Range("D22") "6.905"
$ws.Range("E22").Value = "  -3.98%  "
Set-TextValue $ws.Range("D23") "1.005"
$ws.Range("E23").Value = "  +0.26%  "
Set-TextValue $ws.Range("D24") "156.64"
$ws.Range("E24").Value = "  -0.34%  "
Set-TextValue $ws.Range("D25") "7.982"
$ws.Range("E25").Value = "  -4.91%  "
Set-TextValue $ws.Range("D26") "0.1292"
$ws.Range("E26").Value = "  -2.71%  "
Set-TextValue $ws.Range("D27") "16.54"
$ws.Range("E27").Value = "  -3.73%  "
Set-TextValue $ws.Range("D28") "0.06452"
$ws.Range("E28").Value = "  -6.28%  "
Set-TextValue $ws.Range("D29") "1.422"
$ws.Range("E29").Value = "  -3.54%  "
Set-TextValue $ws.Range("D30") "1.444"
$ws.Range("E30").Value = "  -2.28%  "
Set-TextValue $ws.Range("D31") "3.838"
$ws.Range("E31").Value = "  -2.46%  "
Set-TextValue $ws.Range("D32") "3.800"
$ws.Range("E32").Value = "  -5.31%  "
Set-TextValue $ws.Range("D33") "1.096"
$ws.Range("E33").Value = "  -4.96%  "
Set-TextValue $ws.Range("D34") "1.729"
$ws.Range("E34").Value = "  -0.78%  "
Set-TextValue $ws.Range("D35") "0.6481"
$ws.Range("E35").Value = "  -4.53%  "
$ws.Range("E36").Value = "  -1.36%  "
Set-TextValue $ws.Range("D37") "2.756"
$ws.Range("E37").Value = "  -0.80%  "
Set-TextValue $ws.Range("D38") "0.01756"
$ws.Range("E38").Value = "  -3.52%  "
Set-TextValue $ws.Range("D39") "6.617"
$ws.Range("E39").Value = "  -0.62%  "
Set-TextValue $ws.Range("D40") "1.142.01"
$ws.Range("E40").Value = "  -7.14%  "
Set-TextValue $ws.Range("D41") "0.8941"
$ws.Range("E41").Value = "  -5.41%  "
Set-TextValue $ws.Range("D42") "1.005"
$ws.Range("E42").Value = "  +0.37%  "
Set-TextValue $ws.Range("D43") "2.001.45"
$ws.Range("E43").Value = "  +0.56%  "
Set-TextValue $ws.Range("D44") "100.76"
$ws.Range("E44").Value = "  -0.33%  "
Set-TextValue $ws.Range("D45") "62.29"
$ws.Range("E45").Value = "  -4.27%  "
$ws.Range("E46").Value = "  -3.40%  "
Set-TextValue $ws.Range("D47") "1.624"
$ws.Range("E47").Value = "  -3.70%  "
Set-TextValue $ws.Range("D48") "8.519"
$ws.Range("E48").Value = "  -3.22%  "
Set-TextValue $ws.Range("D49") "0.4542"
$ws.Range("E49").Value = "  -0.73%  "
Set-TextValue $ws.Range("D50") "0.05497"
$ws.Range("E50").Value = "  -2.51%  "
Set-TextValue $ws.Range("D51") "6.392"
$ws.Range("E51").Value = "  -6.85%  "
